$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")

# New rows of data describing "treatment 18" - a calendar month view with
# icon and titration interaction - modeled on the existing "treatment 17"
# block (rows 32:34) directly above.
$comment = "Calendar month view with icon and titration interaction."

$rows = @(
    @{ row = 35; B = 1; F = 300; H = 44593; I = 700;  K = 44614; },
    @{ row = 36; B = 2; F = 500; H = 44621; I = 800;  K = 44632; },
    @{ row = 37; B = 3; F = 300; H = 44652; I = 1000; K = 44666; }
)

foreach ($r in $rows) {
    $i = $r.row
    $ws.Cells.Item($i, 1).Value = 18            # A - treatment_id
    $ws.Cells.Item($i, 2).Value = $r.B          # B - position
    $ws.Cells.Item($i, 3).Value = "calendarIcon" # C - view_type
    $ws.Cells.Item($i, 4).Value = "titration"    # D - interaction
    $ws.Cells.Item($i, 5).Value = "laterAmount"  # E - variable_amount
    $ws.Cells.Item($i, 6).Value = $r.F          # F - amount_earlier
    $ws.Cells.Item($i, 8).Value = $r.H          # H - date_earlier
    $ws.Cells.Item($i, 8).NumberFormat = "m/d/yyyy;@"
    $ws.Cells.Item($i, 9).Value = $r.I          # I - amount_later
    $ws.Cells.Item($i, 11).Value = $r.K         # K - date_later
    $ws.Cells.Item($i, 11).NumberFormat = "m/d/yyyy;@"
    $ws.Cells.Item($i, 12).Value = 1100          # L - max_amount
    $ws.Cells.Item($i, 14).Value = 100           # N - horizontal_pixels
    $ws.Cells.Item($i, 15).Value = 100           # O - vertical_pixels
    $ws.Cells.Item($i, 20).Value = 8             # T - width_in
    $ws.Cells.Item($i, 21).Value = 8             # U - height_in
    $ws.Cells.Item($i, 22).Value = $comment      # V - comment
}

$ws.Range("A20").Select()

$wb.Save()
